# RallyBike/FrameCompare.xlsx - "race 3 2 frames faster"
#
# V4 sheet (sheet1): fill in a previously-blank start-frame (B20),
# insert a new split row "Turbo 180" at row 21, and shave 2 frames off
# the two "race 3" checkpoints that used to be rows 21-22 (now 22-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# B20 was empty; fill in the start frame for this split (diff becomes 686).
$ws.Range("B20").Value = 15514

# Insert a new row for the "Turbo 180" split right after row 20. Excel
# copies the formatting of the row above (styles 17/2/2/3), matching
# rows 19-20.
$ws.Rows(21).Insert() | Out-Null

$ws.Range("A21").Value = "Turbo 180"
$ws.Range("B21").Value = 15905

# The old rows 21 and 22 (now 22 and 23, "Cross finish" / "Black screen"
# for race 3) got 2 frames faster.
$ws.Range("B22").Value = 17673
$ws.Range("B23").Value = 17924

# Match the final cursor position recorded in the saved workbook.
$ws.Range("B24").Select() | Out-Null
